$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BUGS & Issues")

# --- Issue reference codes, column A rows 2-15 ---
$ws.Range("A2").Value = "ISS.001"
$ws.Range("A3").Value = "ISS.002"
$ws.Range("A4").Value = "ISS.003"
$ws.Range("A5").Value = "ISS.004"
$ws.Range("A6").Value = "ISS.005"
$ws.Range("A7").Value = "ISS.006"
$ws.Range("A8").Value = "ISS.007"
$ws.Range("A9").Value = "ISS.008"
$ws.Range("A10").Value = "ISS.009"
$ws.Range("A11").Value = "ISS.010"
$ws.Range("A12").Value = "ISS.011"
$ws.Range("A13").Value = "ISS.012"
$ws.Range("A14").Value = "ISS.013"
$ws.Range("A15").Value = "ISS.014"

# --- Date Raised for the first issue ---
$ws.Range("B2").Value = 44000
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# --- Description for the first issue ---
$ws.Range("D2").Value = "Use of gem 'tzinfo-data' in AWS deployment on Linux servers generates Warning"

# --- Header row ---
$ws.Range("B1").Value = "Date Raised"
$ws.Range("D1").Value = "Description"
$ws.Range("C1").Value = "Status"
$ws.Range("E1").Value = "Assigned to"

# --- Remaining detail for the first issue ---
$ws.Range("C2").Value = "In Progress"
$ws.Range("E2").Value = "Mark C"

# --- Ref header (reuses existing shared string) ---
$ws.Range("A1").Value = "Ref"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 75

# --- Activate the BUGS & Issues sheet and select E3, matching the saved view state ---
$ws.Activate() | Out-Null
$ws.Range("E3").Select() | Out-Null
